# List of Fields in STARS.xlsx -- "Add files via upload"
#
# Adds a large batch of new field rows (22-41: Social History / Bar Code ID /
# Financial sections) plus backfills the "Table Name" (column B) values and a
# handful of missing Required/Demanded/Necessary + Entry Method cells on the
# existing rows (4-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New rows 22-41 : columns A-F first (row by row), matching the order the
#    values were originally typed in (keeps the shared-string table order
#    identical to the source edit).
# ---------------------------------------------------------------------------

$ws.Range("A22").Value = "Veteran"
$ws.Range("B22").Value = "Social History"
$ws.Range("D22").Value = "Yes"
$ws.Range("F22").Value = "Button"

$ws.Range("A23").Value = "Client diasabled"
$ws.Range("B23").Value = "Social History"
$ws.Range("D23").Value = "Yes"
$ws.Range("E23").Value = "Yes"
$ws.Range("F23").Value = "Button"

$ws.Range("A24").Value = "Lives With"
$ws.Range("B24").Value = "Social History"
$ws.Range("D24").Value = "Yes"
$ws.Range("F24").Value = "Drop Down"

$ws.Range("A25").Value = "Does not speak english"
$ws.Range("B25").Value = "Social History"
$ws.Range("F25").Value = "Button"

$ws.Range("A26").Value = "Primary Language"
$ws.Range("B26").Value = "Social History"
$ws.Range("C26").Value = "Yes"
$ws.Range("E26").Value = "Yes"
$ws.Range("F26").Value = "Drop Down"

$ws.Range("A27").Value = "Special needs"
$ws.Range("B27").Value = "Social History"
$ws.Range("E27").Value = "Yes"
$ws.Range("F27").Value = "List Select"

$ws.Range("A28").Value = "Race"
$ws.Range("B28").Value = "Social History"
$ws.Range("D28").Value = "Yes"
$ws.Range("E28").Value = "Yes"
$ws.Range("F28").Value = "Check Mark"

$ws.Range("A29").Value = "Ethnicity"
$ws.Range("B29").Value = "Social History"
$ws.Range("D29").Value = "Yes"
$ws.Range("E29").Value = "Yes"
$ws.Range("F29").Value = "Drop Down"

$ws.Range("A30").Value = "Oxygen dependent"
$ws.Range("B30").Value = "Social History"
$ws.Range("E30").Value = "Yes"
$ws.Range("F30").Value = "Button"

$ws.Range("A31").Value = "Insulin Dependent"
$ws.Range("B31").Value = "Social History"
$ws.Range("E31").Value = "Yes"
$ws.Range("F31").Value = "Button"

$ws.Range("A32").Value = "Dialysis"
$ws.Range("B32").Value = "Social History"
$ws.Range("E32").Value = "Yes"
$ws.Range("F32").Value = "Button"

$ws.Range("A33").Value = "Community Emergeny high risk"
$ws.Range("B33").Value = "Social History"
$ws.Range("F33").Value = "Button"

$ws.Range("A34").Value = "Status"
$ws.Range("B34").Value = "Bar Code ID"
$ws.Range("C34").Value = "Yes"
$ws.Range("E34").Value = "Yes"
$ws.Range("F34").Value = "Drop Down"

$ws.Range("A35").Value = "Bar code ID"
$ws.Range("B35").Value = "Bar Code ID"
$ws.Range("E35").Value = "Yes"
$ws.Range("F35").Value = "Manual Input"

$ws.Range("A36").Value = "Poverty level"
$ws.Range("B36").Value = "Financial"
$ws.Range("D36").Value = "Yes"
$ws.Range("E36").Value = "Yes"
$ws.Range("F36").Value = "Drop Down"

$ws.Range("A37").Value = "Low Income minority"
$ws.Range("B37").Value = "Financial"
$ws.Range("F37").Value = "Button"

$ws.Range("A38").Value = "Monthly Income"
$ws.Range("B38").Value = "Financial"
$ws.Range("D38").Value = "Yes"
$ws.Range("E38").Value = "Yes"
$ws.Range("F38").Value = "Manual Input"

$ws.Range("A39").Value = "# of Household "
$ws.Range("B39").Value = "Financial"
$ws.Range("D39").Value = "Yes"
$ws.Range("F39").Value = "Manual Input"

$ws.Range("A40").Value = "Income Range"
$ws.Range("B40").Value = "Financial"
$ws.Range("D40").Value = "Yes"
$ws.Range("E40").Value = "Yes"
$ws.Range("F40").Value = "Drop Down"

$ws.Range("A41").Value = "Financial comments"
$ws.Range("B41").Value = "Financial"
$ws.Range("F41").Value = "Manual Input"

# ---------------------------------------------------------------------------
# 2) New rows 22-41 : column G "comments" cells, filled in as a second pass
#    (again mirroring the original authoring order). Rows 30-32 were pasted
#    in without the column's default formatting, so explicitly reset those
#    three back to the workbook's Normal style.
# ---------------------------------------------------------------------------

$ws.Range("G24").Value = "Lives With"
$ws.Range("G26").Value = "Primary Language"
$ws.Range("G27").Value = "Special Communication needs"

$ws.Range("G30").Value = "Oxygen dependent"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "Insulin Dependent"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "Dialysis"
$ws.Range("G32").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Backfill column B ("Table Name" = Basic Demographics) on the existing
#    rows 4-19, plus a few previously-missing Required/Demanded/Necessary
#    and Entry Method cells that were filled in at the same time.
# ---------------------------------------------------------------------------

$ws.Range("B4").Value = "Basic Demographics"
$ws.Range("B5").Value = "Basic Demographics"
$ws.Range("B6").Value = "Basic Demographics"
$ws.Range("B7").Value = "Basic Demographics"
$ws.Range("B8").Value = "Basic Demographics"
$ws.Range("B9").Value = "Basic Demographics"
$ws.Range("B10").Value = "Basic Demographics"
$ws.Range("B11").Value = "Basic Demographics"
$ws.Range("B12").Value = "Basic Demographics"
$ws.Range("B13").Value = "Basic Demographics"
$ws.Range("B14").Value = "Basic Demographics"
$ws.Range("B15").Value = "Basic Demographics"
$ws.Range("B16").Value = "Basic Demographics"
$ws.Range("B17").Value = "Basic Demographics"
$ws.Range("B18").Value = "Basic Demographics"
$ws.Range("B19").Value = "Basic Demographics"

$ws.Range("G14").Value = "Lives in an elevator building"
$ws.Range("D15").Value = "Yes"
$ws.Range("E15").Value = "Yes"
$ws.Range("E16").Value = "Yes"
$ws.Range("E17").Value = "Yes"
$ws.Range("E18").Value = "Yes"

# ---------------------------------------------------------------------------
# 4) Restore the selection to where the author left off.
# ---------------------------------------------------------------------------

$ws.Range("H35").Select() | Out-Null
